# Fix some SMT positions for pick and place
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C6: Mid Y -99.0000mm -> -97.0000mm
$ws.Range("C8").Value = "-97.0000mm"

# R2: Mid X 172.0000mm -> 169.0000mm
$ws.Range("B20").Value = "169.0000mm"

# R3: Mid X 174.5000mm -> 171.5000mm ; Mid Y -119.0625mm -> -119.0000mm
$ws.Range("B21").Value = "171.5000mm"
$ws.Range("C21").Value = "-119.0000mm"

# FB2: Mid X 167.0000mm -> 155.5000mm ; Mid Y -119.0625mm -> -111.5000mm
$ws.Range("B29").Value = "155.5000mm"
$ws.Range("C29").Value = "-111.5000mm"

# FB3: Mid X 154.0625mm -> 153.0000mm ; Mid Y -110.5000mm -> -111.5000mm ; Rotation 0 -> 270
$ws.Range("B30").Value = "153.0000mm"
$ws.Range("C30").Value = "-111.5000mm"
$ws.Range("D30").Value = 270
